# Update "想去人数" (column F, numeric) values as published for the gh-pages
# data snapshot generated at commit 456a3b4.
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances) -- unchanged
# Sheet 3 = 本地生活 (Local Life)
# Sheet 4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item(1)   # 展览
$wsLocal = $wb.Worksheets.Item(3)  # 本地生活
$wsAll = $wb.Worksheets.Item(4)    # 全部类型

# --- 展览 (sheet 1) ---
$wsExpo.Cells.Item(2, 6).Value = 39
$wsExpo.Cells.Item(5, 6).Value = 979
$wsExpo.Cells.Item(6, 6).Value = 368
$wsExpo.Cells.Item(8, 6).Value = 560
$wsExpo.Cells.Item(9, 6).Value = 1460
$wsExpo.Cells.Item(11, 6).Value = 1343
$wsExpo.Cells.Item(12, 6).Value = 3005
$wsExpo.Cells.Item(13, 6).Value = 425
$wsExpo.Cells.Item(14, 6).Value = 1624
$wsExpo.Cells.Item(16, 6).Value = 796
$wsExpo.Cells.Item(18, 6).Value = 1389
$wsExpo.Cells.Item(19, 6).Value = 268
$wsExpo.Cells.Item(22, 6).Value = 401
$wsExpo.Cells.Item(23, 6).Value = 2
$wsExpo.Cells.Item(24, 6).Value = 3479
$wsExpo.Cells.Item(25, 6).Value = 686
$wsExpo.Cells.Item(27, 6).Value = 1541

# --- 本地生活 (sheet 3) ---
$wsLocal.Cells.Item(3, 6).Value = 4

# --- 全部类型 (sheet 4) ---
$wsAll.Cells.Item(2, 6).Value = 39
$wsAll.Cells.Item(5, 6).Value = 4
$wsAll.Cells.Item(16, 6).Value = 979
$wsAll.Cells.Item(17, 6).Value = 368
$wsAll.Cells.Item(19, 6).Value = 560
$wsAll.Cells.Item(20, 6).Value = 1460
$wsAll.Cells.Item(22, 6).Value = 1343
$wsAll.Cells.Item(23, 6).Value = 3005
$wsAll.Cells.Item(24, 6).Value = 425
$wsAll.Cells.Item(25, 6).Value = 1624
$wsAll.Cells.Item(27, 6).Value = 796
$wsAll.Cells.Item(29, 6).Value = 1389
$wsAll.Cells.Item(30, 6).Value = 268
$wsAll.Cells.Item(35, 6).Value = 401
$wsAll.Cells.Item(36, 6).Value = 2
$wsAll.Cells.Item(37, 6).Value = 3479
$wsAll.Cells.Item(38, 6).Value = 686
$wsAll.Cells.Item(40, 6).Value = 1541
